$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add a new (auto-numbered) Print_Titles defined name, mirroring the
#    repeated "set print titles" pattern already present in the workbook.
#    Scoping it on the worksheet's Names collection keeps it local to
#    Sheet1 (localSheetId 0), matching the existing entries.
# ---------------------------------------------------------------------------
$ws.Names.Add("_xlnm.Print_Titles_0_0_0_0_0_0_0_0_0", "=Sheet1!`$3:`$4")

# ---------------------------------------------------------------------------
# 2) Row 23: drop the trailing, value-less R23 cell (it only carried a
#    default style, no content) by clearing it outright.
# ---------------------------------------------------------------------------
$ws.Cells.Item(23, 18).ClearContents()

# ---------------------------------------------------------------------------
# 3) Append the new data row (row 26) describing "Hirudoid, Creme".
#    Styling is copied implicitly by re-using the same NumberFormat as the
#    equivalent cells on the row above (row 25) so the generated styles
#    match the workbook's existing per-column formatting (ID numbers,
#    dates, package-size counters, ...).
# ---------------------------------------------------------------------------
$r = 26

$ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(25, 1).NumberFormat
$ws.Cells.Item($r, 1).Value2 = 16105

$ws.Cells.Item($r, 2).Value2 = 1

$ws.Cells.Item($r, 3).Value2 = "Hirudoid, Creme"
$ws.Cells.Item($r, 4).Value2 = "Medinova AG"
$ws.Cells.Item($r, 5).Value2 = "02.08.2."
$ws.Cells.Item($r, 6).Value2 = "C05BA"
$ws.Cells.Item($r, 7).Value2 = "Synthetika human"

$ws.Cells.Item($r, 8).NumberFormat = $ws.Cells.Item(25, 8).NumberFormat
$ws.Cells.Item($r, 8).Value2 = 18872

$ws.Cells.Item($r, 9).NumberFormat = $ws.Cells.Item(25, 9).NumberFormat
$ws.Cells.Item($r, 9).Value2 = 18872

$ws.Cells.Item($r, 10).NumberFormat = $ws.Cells.Item(25, 10).NumberFormat
$ws.Cells.Item($r, 10).Value2 = 43162

$ws.Cells.Item($r, 11).NumberFormat = $ws.Cells.Item(25, 11).NumberFormat
$ws.Cells.Item($r, 11).Value2 = 58

# Packungsgroesse "40" is stored as text (matches e.g. "150", "2x10" sibling
# cells in this column), so force text via a leading quote like a user
# typing '40 would in Excel.
$ws.Cells.Item($r, 12).Value2 = "'40"

$ws.Cells.Item($r, 13).Value2 = "g"
$ws.Cells.Item($r, 14).Value2 = "D"
$ws.Cells.Item($r, 15).Value2 = "heparinoidum (chondroitini polysulfas)"
$ws.Cells.Item($r, 16).Value2 = "heparinoidum (chondroitini polysulfas) 3 mg alcoholes adipis lanae, aromatica, conserv.: E 218, E 216, excipiens ad unguentum pro 1 g."
$ws.Cells.Item($r, 17).Value2 = "Venenmittel für den äusserlichen Gebrauch"

# ---------------------------------------------------------------------------
# 4) Nudge the header logo picture's anchor/size slightly (matches the
#    small reflow caused by the extra row of data above it).
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 270000 / 12700
$shp.Top = 0
$shp.Width = 1777680 / 12700
$shp.Height = 539640 / 12700

# ---------------------------------------------------------------------------
# 5) Leave the selection on the newly added row, like the author would have
#    after finishing the edit.
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Select()
